$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BO2").Value = 0.7000501751899719
$ws.Range("BO3").Value = 0.7252546548843384
$ws.Range("BO4").Value = 0.8770965337753296
$ws.Range("BO5").Value = 0.726647675037384
$ws.Range("BO6").Value = 0.752899169921875
$ws.Range("BO7").Value = 0.7230728268623352
$ws.Range("BO8").Value = 0.7775186896324158
$ws.Range("BO9").Value = 0.7579985857009888
$ws.Range("BO10").Value = 0.729798436164856
$ws.Range("BO11").Value = 0.7400994896888733
$ws.Range("BO12").Value = 0.7421934008598328
$ws.Range("BO13").Value = 0.7053585052490234
$ws.Range("BO14").Value = 0.7672739624977112
$ws.Range("BO15").Value = 0.6614869236946106
$ws.Range("BO16").Value = 0.7542528510093689
$ws.Range("BO17").Value = 0.7394630908966064
$ws.Range("BO18").Value = 0.7813862562179565
$ws.Range("BO19").Value = 0.345298707485199
$ws.Range("BO20").Value = 0.7774925231933594
$ws.Range("BO21").Value = 0.4294066727161407
$ws.Range("BO22").Value = 0.4785134196281433
$ws.Range("BO23").Value = 0.7571130990982056
$ws.Range("BO24").Value = 0.6992326378822327
$ws.Range("BO25").Value = 0.6605466604232788
$ws.Range("BO26").Value = 0.71180659532547
$ws.Range("BO27").Value = 0.6626912355422974
$ws.Range("BO28").Value = 0.6601747274398804
$ws.Range("BO29").Value = 0.6637101173400879
$ws.Range("BO30").Value = 0.7701870799064636
$ws.Range("BO31").Value = 0.7080376744270325
$ws.Range("BO32").Value = 0.5439549088478088
$ws.Range("BO33").Value = 0.6939806938171387
$ws.Range("BO34").Value = 0.5539235472679138
$ws.Range("BO35").Value = 0.7407917380332947
$ws.Range("BO36").Value = 0.5644567012786865
$ws.Range("BO37").Value = 0.8088090419769287
$ws.Range("BO38").Value = 0.6603132486343384
$ws.Range("BO39").Value = 0.6831006407737732
$ws.Range("BO40").Value = 0.6939290761947632
$ws.Range("BO41").Value = 0.5631906390190125
$ws.Range("BO42").Value = 0.49684077501297
$ws.Range("BO43").Value = 0.5076735019683838
$ws.Range("BO44").Value = 0.7261407375335693
$ws.Range("BO45").Value = 0.5679957866668701
$ws.Range("BO46").Value = 0.457265168428421
$ws.Range("BO47").Value = 0.7280606031417847
$ws.Range("BO48").Value = 0.7098612785339355
$ws.Range("BO49").Value = 0.7323701977729797
$ws.Range("BO50").Value = 0.7484275102615356
$ws.Range("BO51").Value = 0.6802716851234436
$ws.Range("BO52").Value = 0.322256863117218
$ws.Range("BO53").Value = 0.2477652579545975
$ws.Range("BO54").Value = 0.005304196383804083
$ws.Range("BO55").Value = 0.7294290661811829
$ws.Range("BO56").Value = 0.8103954195976257
$ws.Range("BO57").Value = 0.7873294353485107
$ws.Range("BO58").Value = 0.02227035723626614
$ws.Range("BO59").Value = 0.7275230288505554
$ws.Range("BO60").Value = 0.7685961723327637
$ws.Range("BO61").Value = 0.7116091251373291
$ws.Range("BO62").Value = 0.7875267267227173
$ws.Range("BO63").Value = 0.6278612613677979
$ws.Range("BO64").Value = 0.6266006827354431
$ws.Range("BO65").Value = 0.7793653011322021
$ws.Range("BO66").Value = 0.2175709456205368
$ws.Range("BO67").Value = 0.5856631398200989
$ws.Range("BO68").Value = 0.2039490342140198
$ws.Range("BO69").Value = 0.698826014995575
$ws.Range("BO70").Value = 0.7884714603424072
$ws.Range("BO71").Value = 0.157875120639801
$ws.Range("BO72").Value = 0.7256271243095398
$ws.Range("BO73").Value = 0.7218974232673645
$ws.Range("BO74").Value = 0.7109302878379822
$ws.Range("BO75").Value = 0.7295388579368591
$ws.Range("BO76").Value = 0.6587293744087219
